$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 650, shifting existing rows 650-727 down to 651-728.
$ws.Rows.Item(650).Insert()

# Populate the newly inserted row 650 with the new record's data.
$ws.Range("A650").Value = 4
$ws.Range("B650").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C650").Value = 'Los Lagos'
$ws.Range("D650").Value = 45142
$ws.Range("E650").Value = 10
$ws.Range("F650").Value = 100114001
$ws.Range("G650").Value = 'Papa'
$ws.Range("H650").Value = 'Patagonia'
$ws.Range("I650").Value = '1a (guarda)'
$ws.Range("J650").Value = 600
$ws.Range("K650").Value = 18000
$ws.Range("L650").Value = 19000
$ws.Range("M650").Value = 18500
$ws.Range("N650").Value = '$/saco 25 kilos'
$ws.Range("O650").Value = 'Provincia de Llanquihue'
$ws.Range("P650").Value = 740
$ws.Range("Q650").Value = 25
$ws.Range("R650").Value = 'Hortaliza'

# Match the number format/style used by the other date cells in column D.
$ws.Range("D650").NumberFormat = $ws.Range("D651").NumberFormat
